$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 4
$ws.Range("H2").Value = 3.75
$ws.Range("I2").Value = 1.8
$ws.Range("J2").Value = 4.5
$ws.Range("X2").Value = 21
$ws.Range("Y2").Value = 13
$ws.Range("AI2").Value = 9
$ws.Range("AY2").Value = 21
$ws.Range("AZ2").Value = 34
$ws.Range("I3").Value = 11
$ws.Range("J3").Value = 1.91
$ws.Range("K3").Value = 2.3
$ws.Range("M3").Value = 1.05
$ws.Range("N3").Value = 11
$ws.Range("O3").Value = 1.3
$ws.Range("P3").Value = 3.4
$ws.Range("Q3").Value = 1.98
$ws.Range("R3").Value = 1.83
$ws.Range("U3").Value = 2.38
$ws.Range("V3").Value = 1.53
$ws.Range("W3").Value = 5.5
$ws.Range("AC3").Value = 9
$ws.Range("AE3").Value = 26
$ws.Range("BA3").Value = 301
$ws.Range("G4").Value = 1.9
$ws.Range("H4").Value = 3.6
$ws.Range("I4").Value = 4
$ws.Range("J4").Value = 2.5
$ws.Range("L4").Value = 4.5
$ws.Range("O4").Value = 1.3
$ws.Range("P4").Value = 3.4
$ws.Range("U4").Value = 1.83
$ws.Range("V4").Value = 1.83
$ws.Range("W4").Value = 7
$ws.Range("X4").Value = 9
$ws.Range("Y4").Value = 8.5
$ws.Range("Z4").Value = 15
$ws.Range("AA4").Value = 15
$ws.Range("AI4").Value = 21
$ws.Range("AJ4").Value = 13
$ws.Range("AL4").Value = 34
$ws.Range("AM4").Value = 41
$ws.Range("AO4").Value = 10
$ws.Range("AQ4").Value = 34
$ws.Range("AW4").Value = 6
$ws.Range("AX4").Value = 21
$ws.Range("AZ4").Value = 81
$ws.Range("BA4").Value = 101
$ws.Range("G6").Value = 2.25
$ws.Range("I6").Value = 3.1
$ws.Range("Q6").Value = 1.9
$ws.Range("R6").Value = 1.95
$ws.Range("S6").Value = 1.4
$ws.Range("T6").Value = 2.75
$ws.Range("AB6").Value = 26
$ws.Range("AG6").Value = 201
$ws.Range("AH6").Value = 10
$ws.Range("AI6").Value = 15
$ws.Range("AT6").Value = 2.75
$ws.Range("I7").Value = 2.45
$ws.Range("J7").Value = 4
$ws.Range("L7").Value = 3.25
$ws.Range("W7").Value = 7.5
$ws.Range("AC7").Value = 6.5
$ws.Range("AI7").Value = 10
$ws.Range("BB7").Value = 251
$ws.Range("BD8").Value = 151
$ws.Range("N9").Value = 8.5
$ws.Range("AN9").Value = 4.33
$ws.Range("G12").Value = 3.05
$ws.Range("I12").Value = 2.25
$ws.Range("J12").Value = 3.55
$ws.Range("K12").Value = 2.07
$ws.Range("L12").Value = 2.8
$ws.Range("O12").Value = 1.34
$ws.Range("W12").Value = 8.75
$ws.Range("X12").Value = 15.5
$ws.Range("Y12").Value = 11
$ws.Range("Z12").Value = 40
$ws.Range("AA12").Value = 28
$ws.Range("AB12").Value = 37
$ws.Range("AE12").Value = 14
$ws.Range("AH12").Value = 7.3
$ws.Range("AI12").Value = 10.75
$ws.Range("AJ12").Value = 9
$ws.Range("AK12").Value = 23
$ws.Range("AL12").Value = 19
$ws.Range("AM12").Value = 30
$ws.Range("AN12").Value = 5
$ws.Range("AO12").Value = 16.5
$ws.Range("AP12").Value = 23
$ws.Range("AQ12").Value = 80
$ws.Range("AR12").Value = 110
$ws.Range("AT12").Value = 2.55
$ws.Range("AW12").Value = 4.15
$ws.Range("AX12").Value = 11.5
$ws.Range("AY12").Value = 19
$ws.Range("BA12").Value = 75
